# Applies the weekly update described by the commit message "Fruta / hortaliza, semanal":
# a new observation is inserted as row 31 (pushing all subsequent rows down by one),
# leaving the rest of the historical rows untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 31, shifting rows 31:106 down to 32:107
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new data point
$ws.Cells.Item(31, 1).Value = 7
$ws.Cells.Item(31, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(31, 3).Value = "Ñuble"
$ws.Cells.Item(31, 4).Value = 44925
$ws.Cells.Item(31, 5).Value = 16
$ws.Cells.Item(31, 6).Value = 100112030
$ws.Cells.Item(31, 7).Value = "Poroto granado"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 60
$ws.Cells.Item(31, 11).Value = 32000
$ws.Cells.Item(31, 12).Value = 32000
$ws.Cells.Item(31, 13).Value = 32000
$ws.Cells.Item(31, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(31, 15).Value = "Región del Maule"
$ws.Cells.Item(31, 16).Value = 1280
$ws.Cells.Item(31, 17).Value = 25
$ws.Cells.Item(31, 18).Value = "Hortaliza"
